# Weekly update: insert a new week's worth of data (2 rows) at the top of
# the Brócoli / Vega Monumental Concepción data block, pushing the existing
# rows down by two. Dimension grows from A1:R292 to A1:R294.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows at 191:192 (formatting, e.g. the date style on
# column D, is inherited from the surrounding rows automatically).
$ws.Range("A191:A192").EntireRow.Insert()

# New row 191 - "Primera" quality entry for 2022-06-07 (serial 44719)
$ws.Cells.Item(191, 1).Value = 11
$ws.Cells.Item(191, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(191, 3).Value = "Bíobío"
$ws.Cells.Item(191, 4).Value = 44719
$ws.Cells.Item(191, 5).Value = 8
$ws.Cells.Item(191, 6).Value = 100112023
$ws.Cells.Item(191, 7).Value = "Brócoli"
$ws.Cells.Item(191, 8).Value = "Sin especificar"
$ws.Cells.Item(191, 9).Value = "Primera"
$ws.Cells.Item(191, 10).Value = 2000
$ws.Cells.Item(191, 11).Value = 700
$ws.Cells.Item(191, 12).Value = 800
$ws.Cells.Item(191, 13).Value = 750
$ws.Cells.Item(191, 14).Value = "`$/unidad"
$ws.Cells.Item(191, 15).Value = "Región Metropolitana"
$ws.Cells.Item(191, 16).Value = 750
$ws.Cells.Item(191, 17).Value = 1
$ws.Cells.Item(191, 18).Value = "Hortaliza"

# New row 192 - "Segunda" quality entry for 2022-06-07 (serial 44719)
$ws.Cells.Item(192, 1).Value = 11
$ws.Cells.Item(192, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(192, 3).Value = "Bíobío"
$ws.Cells.Item(192, 4).Value = 44719
$ws.Cells.Item(192, 5).Value = 8
$ws.Cells.Item(192, 6).Value = 100112023
$ws.Cells.Item(192, 7).Value = "Brócoli"
$ws.Cells.Item(192, 8).Value = "Sin especificar"
$ws.Cells.Item(192, 9).Value = "Segunda"
$ws.Cells.Item(192, 10).Value = 1000
$ws.Cells.Item(192, 11).Value = 600
$ws.Cells.Item(192, 12).Value = 600
$ws.Cells.Item(192, 13).Value = 600
$ws.Cells.Item(192, 14).Value = "`$/unidad"
$ws.Cells.Item(192, 15).Value = "Región Metropolitana"
$ws.Cells.Item(192, 16).Value = 600
$ws.Cells.Item(192, 17).Value = 1
$ws.Cells.Item(192, 18).Value = "Hortaliza"
